$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.869.27'
$ws.Range('E2').Value = '  -0.38%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.892.85'
$ws.Range('E3').Value = '  -0.56%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.000'
$ws.Range('E4').Value = '  -0.54%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7768'
$ws.Range('E5').Value = '  -2.04%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '244.72'
$ws.Range('E6').Value = '  +0.54%  '

$ws.Range('E7').Value = '  -0.48%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3147'
$ws.Range('E8').Value = '  -1.84%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07368'
$ws.Range('E9').Value = '  +3.70%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '25.37'
$ws.Range('E10').Value = '  -3.65%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08143'
$ws.Range('E11').Value = '  +0.97%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.7678'
$ws.Range('E12').Value = '  -0.69%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.503'
$ws.Range('E13').Value = '  +3.23%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.907.47'
$ws.Range('E14').Value = '  -0.65%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '92.41'
$ws.Range('E15').Value = '  -0.43%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.220'
$ws.Range('E16').Value = '  +4.55%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '29.878.69'
$ws.Range('E17').Value = '  -0.45%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '13.98'
$ws.Range('E18').Value = '  +0.30%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '245.66'
$ws.Range('E19').Value = '  +0.06%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007852'
$ws.Range('E20').Value = '  +1.51%  '

$ws.Range('B21').Value = 'Chainlink'
$ws.Range('C21').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '8.177'
$ws.Range('E21').Value = '  +0.59%  '

$ws.Range('B22').Value = 'Dai'
$ws.Range('C22').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.001'
$ws.Range('E22').Value = '  -0.40%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.139.30'
$ws.Range('E23').Value = '  -1.06%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.000'
$ws.Range('E24').Value = '  -0.58%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1572'
$ws.Range('E25').Value = '  -2.46%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.443'
$ws.Range('E26').Value = '  +1.03%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '162.79'
$ws.Range('E27').Value = '  -1.78%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.84'
$ws.Range('E28').Value = '  +0.41%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.047'
$ws.Range('E29').Value = '  -2.66%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.450'
$ws.Range('E30').Value = '  +5.31%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.549'
$ws.Range('E31').Value = '  +0.34%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.531'
$ws.Range('E32').Value = '  +0.56%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05629'
$ws.Range('E33').Value = '  -1.34%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.099'
$ws.Range('E34').Value = '  +0.23%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.251'
$ws.Range('E35').Value = '  -1.59%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7665'
$ws.Range('E36').Value = '  +3.92%  '

$ws.Range('E37').Value = '  +0.50%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.645'
$ws.Range('E38').Value = '  -2.51%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01929'
$ws.Range('E39').Value = '  -0.38%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.788'
$ws.Range('E40').Value = '  +0.10%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.164.65'
$ws.Range('E41').Value = '  +13.15%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '74.45'
$ws.Range('E42').Value = '  +2.67%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.4470'
$ws.Range('E43').Value = '  +0.28%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.990'
$ws.Range('E44').Value = '  +0.62%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.8525'
$ws.Range('E45').Value = '  +0.70%  '

$ws.Range('B46').Value = 'RenderToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.908'
$ws.Range('E46').Value = '  +1.01%  '

$ws.Range('B47').Value = 'PaxDollar'
$ws.Range('C47').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.001'
$ws.Range('E47').Value = '  -0.60%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '102.29'
$ws.Range('E48').Value = '  -0.01%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.922'
$ws.Range('E49').Value = '  +0.80%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '3.087'
$ws.Range('E50').Value = '  +1.30%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.552'
$ws.Range('E51').Value = '  +0.54%  '
